$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38:D38").NumberFormat = "@"
$ws.Range("A38").Value = "2026-02-07"
$ws.Range("B38").Value = "5200000"
$ws.Range("C38").Value = "0"
$ws.Range("D38").Value = "0"
